# Applies the Learning_Python_syllabus.docx content update:
# - refreshed 'Generated on' timestamp
# - reworked week-by-week topics/subtopics/activities text
# Each target run of text is located by its known, pre-computed
# character offsets in $d.Content (verified against the expected
# old text before writing) and replaced via Range.Text so Word's
# AutoCorrect (e.g. smart quotes) cannot mangle the new text.
$d = $word.ActiveDocument
$ops = @(
    ,@(4196, 4250, '* Activity:  Mini project presentation and final quiz.', '* Activity: Mini-Project 2: Develop a data parser or another application of your choice; Final Project Presentations.')
    ,@(4104, 4195, '* Subtopics: Finishing the mini project, code review and debugging, course summary and Q&A.', '* Subtopics:  Design, implementation, testing, presentation.')
    ,@(4048, 4103, '* Main Topic: Mini Project Completion and Course Review', '* Main Topic: Project Development & Review')
    ,@(4004, 4047, '**Week 15: Mini Project - Part 2 & Review**', '**Week 15: Mini-Project 2: Data Parser or Other Application**')
    ,@(3952, 4001, '* Activity:  Project planning and initial coding.', '* Activity:  Mini-Project 1:  Develop a calculator or quiz application')
    ,@(3851, 3951, '* Subtopics: Planning and designing a mini project (e.g., calculator, quiz app, simple data parser).', '* Subtopics:  Design, implementation, testing.')
    ,@(3812, 3850, '* Main Topic: Mini Project Development', '* Main Topic: Project Development')
    ,@(3777, 3811, '**Week 14: Mini Project - Part 1**', '**Week 14: Mini-Project 1: Calculator or Quiz App**')
    ,@(3726, 3774, '* Activity:  Lab:  Parsing data from a CSV file.', '* Activity: File I/O Exercises & "Data Parser" Lab (CSV or JSON)')
    ,@(3654, 3725, '* Subtopics: Reading and writing files, working with CSV and JSON data.', '* Subtopics: Reading and writing files, working with CSV files, working with JSON files.')
    ,@(3564, 3603, '**Week 13: File I/O and Data Handling**', '**Week 13: File I/O and Data Serialization**')
    ,@(3470, 3561, '* Activity:  Lab:  Building a program using inheritance (e.g., different types of animals).', '* Activity: OOP Principles Exercises & "Shape Inheritance" Lab')
    ,@(3415, 3469, '* Subtopics: Inheritance, polymorphism, encapsulation.', '* Subtopics: Inheritance (creating subclasses), encapsulation (data hiding), polymorphism.')
    ,@(3379, 3414, '* Main Topic: Advanced OOP Concepts', '* Main Topic: OOP Principles')
    ,@(3327, 3378, '**Week 12: Object-Oriented Programming (OOP) - II**', '**Week 12: OOP: Inheritance and Encapsulation**')
    ,@(3259, 3324, '* Activity:  Lab:  Creating a simple class (e.g., a `Dog` class).', '* Activity: OOP Basics Exercises & "Simple Animal Class" Lab')
    ,@(3174, 3258, '* Subtopics: Classes and objects, attributes and methods, constructors (`__init__`).', '* Subtopics: Defining classes, creating objects, attributes, methods, constructors (`__init__`).')
    ,@(3140, 3173, '* Main Topic: Introduction to OOP', '* Main Topic: Classes and Objects')
    ,@(3089, 3139, '**Week 11: Object-Oriented Programming (OOP) - I**', '**Week 11: Introduction to Object-Oriented Programming (OOP)**')
    ,@(3022, 3086, '* Activity:  Lab: Building a program with robust error handling.', '* Activity: Exception Handling Exercises & "Robust File Reader" Lab')
    ,@(2938, 3021, '* Subtopics: `try`, `except`, `finally` blocks, handling different exception types.', '* Subtopics: try, except, else, finally blocks, common exceptions (TypeError, ValueError, FileNotFoundError).')
    ,@(2780, 2859, '* Activity:  Lab: Using a library like `requests` to fetch data from a website.', '* Activity: Module Usage Exercises & "Dice Rolling Simulator" Lab (using random module)')
    ,@(2645, 2779, '* Subtopics: Importing modules (math, random, etc.), creating custom modules, installing packages using pip, using external libraries.', '* Subtopics: Importing modules, using built-in modules (math, random, os), installing packages with pip, using external packages.')
    ,@(2527, 2568, '* Activity:  Quiz on functions and scope.', '* Activity: Function Exercises & "Simple Math Library" Lab')
    ,@(2425, 2526, '* Subtopics: Defining functions, function parameters and arguments, return values, scope, docstrings.', '* Subtopics: Defining functions, parameters and arguments, return values, scope, docstrings.')
    ,@(2386, 2424, '* Main Topic: Functions and Modularity', '* Main Topic: Functions and Modular Programming')
    ,@(2296, 2361, '* Activity: Lab:  Looping exercises (e.g., factorial calculator).', '* Activity: Dictionary Exercises & "Contact Book" Lab')
    ,@(2168, 2295, '* Subtopics: `for` loops, `while` loops, iterating through lists, dictionaries, and strings, `break` and `continue` statements.', '* Subtopics: Dictionary creation, accessing values, methods (keys(), values(), items()), iterating through dictionaries.')
    ,@(2148, 2167, '* Main Topic: Loops', '* Main Topic: Dictionaries and Data Structures')
    ,@(2115, 2147, '**Week 7: Control Flow - Loops**', '**Week 7: Dictionaries**')
    ,@(2026, 2112, '* Activity:  Lab:  Building a simple decision-making program (e.g., grade calculator).', '* Activity: Collection Manipulation Exercises & "Student Grade Manager" Lab (using lists)')
    ,@(1925, 2025, '* Subtopics: `if`, `elif`, `else` statements, nested conditionals, logical operators (and, or, not).', '* Subtopics: Lists (creation, manipulation, methods), Tuples (immutability), Sets (unique elements, set operations).')
    ,@(1888, 1924, '* Main Topic: Conditional Statements', '* Main Topic: Working with Collections')
    ,@(1838, 1887, '**Week 6: Control Flow - Conditional Statements**', '**Week 6: Lists, Tuples, and Sets**')
    ,@(1772, 1835, '* Activity: Case study: analyzing a dataset using dictionaries.', '* Activity: Looping Exercises & "Fibonacci Sequence Generator" Lab')
    ,@(1601, 1771, '* Subtopics: Creating, accessing, modifying dictionaries and sets, dictionary methods (keys(), values(), items(), etc.), set operations (union, intersection, difference).', '* Subtopics: for loops, while loops, break and continue statements, Looping through sequences (strings, lists).')
    ,@(1565, 1600, '* Main Topic: Dictionaries and Sets', '* Main Topic: Loops and Iterations')
    ,@(1531, 1564, '**Week 5: Dictionaries and Sets**', '**Week 5: Control Flow: Loops**')
    ,@(1473, 1528, '* Activity: Lab: List and tuple manipulation exercises.', '* Activity: Conditional Logic Exercises & "Number Guessing Game" Lab')
    ,@(1340, 1472, '* Subtopics: Creating, accessing, modifying lists and tuples, list methods (append(), insert(), remove(), etc.), tuple immutability.', '* Subtopics: if, elif, else statements, Nested conditional statements, Boolean logic (and, or, not).')
    ,@(1309, 1339, '* Main Topic: Lists and Tuples', '* Main Topic: Conditional Statements')
    ,@(1280, 1308, '**Week 4: Lists and Tuples**', '**Week 4: Control Flow: Conditional Statements**')
    ,@(1201, 1277, '* Activity:  Lab:  String manipulation exercises (e.g., palindrome checker).', '* Activity: String Manipulation Quiz & Case Study: Analyzing Text Data')
    ,@(1085, 1200, '* Subtopics: String slicing, indexing, concatenation, methods (upper(), lower(), split(), etc.), string formatting.', '* Subtopics: String slicing, String methods (upper(), lower(), split(), join()), String formatting, Comparison operators.')
    ,@(1006, 1049, '**Week 3: Strings and String Manipulation**', '**Week 3: String Manipulation and Operators**')
    ,@(957, 1003, '* Activity:  Quiz on data types and operators.', '* Activity:  Data Type Exercises & Simple Calculator Lab')
    ,@(855, 956, '* Subtopics: Integers, Floats, Strings, Booleans, Operators (+, -, *, /, //, %, **), Type conversion.', '* Subtopics:  Variables and assignment, Integers, Floats, Strings, Booleans, Basic Arithmetic Operators, Type Conversion.')
    ,@(806, 854, '* Main Topic: Data types and Operators in Python', '* Main Topic: Variables, Data Types, Operators')
    ,@(769, 805, '**Week 2: Data Types and Operators**', '**Week 2: Basic Syntax and Data Types**')
    ,@(704, 766, '* Activity:  Hands-on exercise: printing different data types.', '* Activity:  Setup Quiz & First Python Program')
    ,@(554, 703, '* Subtopics: What is Python?, Setting up your environment (installation, IDE), First Python program ("Hello, world!"), basic output (print function).', '* Subtopics: What is Python?, Why Python?, Installing Python, Choosing an IDE (VS Code, PyCharm, Thonny),  Hello World program.')
    ,@(500, 553, '* Main Topic:  Introduction to Programming and Python', '* Main Topic: Introduction to Python, Setting up Development Environment (IDE/Text Editor)')
    ,@(457, 499, '**Week 1: Introduction to Python & Setup**', '**Week 1: Introduction to Python & Setting up your environment**')
    ,@(101, 454, '**Course Objectives:** Understand Python Syntax & Structure; Write Python Programs; Use Data Types & Operators Effectively; Implement Control Structures; Modularize Code with Functions & Modules; Handle Errors & Exceptions; Apply Object-Oriented Programming (OOP); Work with Libraries & Packages; Manipulate Files & Data; Build Real-world Mini Projects.', '**Course Objectives:**  Understand Python Syntax & Structure; Write Python Programs; Use Data Types & Operators Effectively; Implement Control Structures; Modularize Code with Functions & Modules; Handle Errors & Exceptions; Apply Object-Oriented Programming (OOP); Work with Libraries & Packages; Manipulate Files & Data; Build Real-world Mini Projects.')
    ,@(62, 99, '## Learning Python - 15-Week Syllabus', '**Learning Python - 15-Week Syllabus**')
    ,@(27, 60, 'Generated on: 2025-06-16 14:02:41', 'Generated on: 2025-06-18 13:24:38')
)

$applied = 0
$skipped = 0
foreach ($op in $ops) {
    $start = $op[0]
    $end = $op[1]
    $expectedOld = $op[2]
    $newText = $op[3]
    $rng = $d.Range($start, $end)
    if ($rng.Text -eq $expectedOld) {
        $rng.Text = $newText
        $applied = $applied + 1
    } else {
        # Fallback: the text must have shifted (e.g. a prior edit moved things)
        # around) -- locate it with Find instead of trusting the stale offset.
        $search = $d.Content
        $ok = $search.Find.Execute($expectedOld, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
        if ($ok) {
            $applied = $applied + 1
        } else {
            $skipped = $skipped + 1
            Write-Output "MISSED: $expectedOld"
        }
    }
}
Write-Output "Applied $applied of $($ops.Count) replacements ($skipped skipped)."
